$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 8816.5
$ws.Range("I131").Value = 673.7778
$ws.Range("J131").Value = 19285.715
$ws.Range("K131").Value = 2021.3334
$ws.Range("L131").Value = 57857.145
$ws.Range("M131").Value = 3018.6666
$ws.Range("N131").Value = -67937.145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 182667.28
$ws.Range("I132").Value = 4512.5
$ws.Range("J132").Value = 503345.9
$ws.Range("K132").Value = 13537.5
$ws.Range("L132").Value = 1510037.7
$ws.Range("M132").Value = -11007.5
$ws.Range("N132").Value = -1515097.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 13159199
$ws.Range("I135").Value = 445
$ws.Range("J135").Value = 23811524
$ws.Range("K135").Value = 4005
$ws.Range("L135").Value = 214303716
$ws.Range("M135").Value = -1470
$ws.Range("N135").Value = -214308786

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 28442.676
$ws.Range("I137").Value = 48481.332
$ws.Range("J137").Value = 6294.684
$ws.Range("K137").Value = 145443.996
$ws.Range("L137").Value = 18884.052
$ws.Range("M137").Value = -142893.996
$ws.Range("N137").Value = -23984.052

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 17171.387
$ws.Range("I141").Value = 1019.7857
$ws.Range("J141").Value = 62395.867
$ws.Range("K141").Value = 3059.3571
$ws.Range("L141").Value = 187187.601
$ws.Range("M141").Value = 2120.6429
$ws.Range("N141").Value = -197547.601

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4181.1
$ws.Range("I32").Value = 4085.8691
$ws.Range("J32").Value = 4681.0625
$ws.Range("K32").Value = 4085.8691
$ws.Range("L32").Value = 4681.0625
$ws.Range("M32").Value = -3798.8691
$ws.Range("N32").Value = -5255.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2388450
$ws.Range("I132").Value = 2837411.5
$ws.Range("J132").Value = 919122.0600000001
$ws.Range("K132").Value = 8512234.5
$ws.Range("L132").Value = 2757366.18
$ws.Range("M132").Value = -8509704.5
$ws.Range("N132").Value = -2762426.18

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 238.4
$ws.Range("I64").Value = 147.33333
$ws.Range("J64").Value = 277.42856
$ws.Range("K64").Value = 147.33333
$ws.Range("L64").Value = 277.42856
$ws.Range("M64").Value = 77.66667000000001
$ws.Range("N64").Value = -727.4285600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 238.4
$ws.Range("I67").Value = 147.33333
$ws.Range("J67").Value = 277.42856
$ws.Range("K67").Value = 147.33333
$ws.Range("L67").Value = 277.42856
$ws.Range("M67").Value = 632.6666700000001
$ws.Range("N67").Value = -1837.42856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 952.375
$ws.Range("I16").Value = 952.375
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 952.375
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -665.375
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 369.4762
$ws.Range("J107").Value = 742.4
$ws.Range("L107").Value = 742.4
$ws.Range("N107").Value = -4582.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 952.375
$ws.Range("I113").Value = 952.375
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 952.375
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1217.625
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 19926
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 19926
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 19926
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -27184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1307.6666
$ws.Range("I122").Value = 1307.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3922.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1472.9998
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1500.4667
$ws.Range("I132").Value = 1464.2609
$ws.Range("J132").Value = 1619.4286
$ws.Range("K132").Value = 4392.7827
$ws.Range("L132").Value = 4858.2858
$ws.Range("M132").Value = -1862.7827
$ws.Range("N132").Value = -9918.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1119.9697
$ws.Range("I134").Value = 951.0417
$ws.Range("J134").Value = 1570.4445
$ws.Range("K134").Value = 2853.1251
$ws.Range("L134").Value = 4711.333500000001
$ws.Range("M134").Value = -318.1251000000002
$ws.Range("N134").Value = -9781.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 5233
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 6877.3335
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 20632.0005
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -21624.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 78125544
$ws.Range("J131").Value = 156250690
$ws.Range("L131").Value = 468752070
$ws.Range("N131").Value = -468762150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1012.0625
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 1107.75
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 1107.75
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -1697.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1012.0625
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 1107.75
$ws.Range("K27").Value = 725
$ws.Range("L27").Value = 1107.75
$ws.Range("M27").Value = -618
$ws.Range("N27").Value = -1321.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 61724
$ws.Range("I40").Value = 1935.5
$ws.Range("J40").Value = 114869.336
$ws.Range("K40").Value = 1935.5
$ws.Range("L40").Value = 114869.336
$ws.Range("M40").Value = -1799.5
$ws.Range("N40").Value = -115141.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 467119.66
$ws.Range("I132").Value = 158152.84
$ws.Range("J132").Value = 776086.4399999999
$ws.Range("K132").Value = 474458.52
$ws.Range("L132").Value = 2328259.32
$ws.Range("M132").Value = -471928.52
$ws.Range("N132").Value = -2333319.32

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 205255.89
$ws.Range("I136").Value = 294866.28
$ws.Range("J136").Value = 2139
$ws.Range("K136").Value = 884598.8400000001
$ws.Range("L136").Value = 6417
$ws.Range("M136").Value = -882048.8400000001
$ws.Range("N136").Value = -11517

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3407.2727
$ws.Range("I62").Value = 3275.5557
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3275.5557
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2651.5557
$ws.Range("N62").Value = -5248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3407.2727
$ws.Range("I65").Value = 3275.5557
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 16377.7785
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -13257.7785
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1147.25
$ws.Range("I81").Value = 768.5333000000001
$ws.Range("J81").Value = 1778.4445
$ws.Range("K81").Value = 1537.0666
$ws.Range("L81").Value = 3556.889
$ws.Range("M81").Value = -476.0666000000001
$ws.Range("N81").Value = -5678.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1147.25
$ws.Range("I84").Value = 768.5333000000001
$ws.Range("J84").Value = 1778.4445
$ws.Range("K84").Value = 7685.333000000001
$ws.Range("L84").Value = 17784.445
$ws.Range("M84").Value = -2381.333000000001
$ws.Range("N84").Value = -28392.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4042.0293
$ws.Range("I132").Value = 1125.2106
$ws.Range("J132").Value = 7736.6665
$ws.Range("K132").Value = 3375.6318
$ws.Range("L132").Value = 23209.9995
$ws.Range("M132").Value = -845.6318000000001
$ws.Range("N132").Value = -28269.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 840208.3
$ws.Range("I136").Value = 1021106
$ws.Range("J136").Value = 371214.22
$ws.Range("K136").Value = 3063318
$ws.Range("L136").Value = 1113642.66
$ws.Range("M136").Value = -3060768
$ws.Range("N136").Value = -1118742.66
